$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing bug statuses from PENDIENTE to CORREGIDO ---
$ws.Range("F5").Value = "CORREGIDO"
$ws.Range("F31").Value = "CORREGIDO"
$ws.Range("F34").Value = "CORREGIDO"
$ws.Range("F35").Value = "CORREGIDO"
$ws.Range("F36").Value = "CORREGIDO"

# --- Fill in the missing responsable for bug #30 ---
$ws.Range("D32").Value = "Pau"

# --- Add new bug row (#35): "Buscador en los Combos" ---
$ws.Range("B36").Copy()
$ws.Range("B37").PasteSpecial(-4122)
$ws.Range("C36").Copy()
$ws.Range("C37").PasteSpecial(-4122)
$ws.Range("F36").Copy()
$ws.Range("F37").PasteSpecial(-4122)

$ws.Range("B37").Value = "Buscador en los Combos"
$ws.Range("C37").Value = "Implementar buscador en los cambios"
$ws.Range("D37").Value = "-"
$ws.Range("E37").Value = "todos"
$ws.Range("F37").Value = "PENDIENTE"

# --- Update the active selection ---
$ws.Range("J30").Select()
